$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete data rows (8-14); this also drops the shared
# strings / styles that were only referenced there.
$ws.Range("A8:A14").EntireRow.Delete()

# Update the report title/header text for the new period.
$ws.Range("A2").Value = "Rekap Penerimaan & AGING 2017/2018 Genap"

# Match the author's final selection in the saved workbook.
$ws.Range("B9").Select()
